$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet: "C. Batumi" -> "C.Batumi"
$ws.Name = "C.Batumi"

# Copy formatting from column J into column K for rows 3-6 (xlPasteFormats = -4122)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)

# Set the new values for year 2023
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1577.4
$ws.Range("K5").Value = 1265.7
$ws.Range("K6").Value = 1824.7

$excel.CutCopyMode = 0
